$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 36.94436433333333
$ws.Range("H2").Value = 110.833093
$ws.Range("I2").Value = 0.8328964975864823
$ws.Range("J2").Value = 0.8328964975864824
$ws.Range("M2").Value = 186.2750726666667
$ws.Range("N2").Value = 558.825218
$ws.Range("O2").Value = 0.9729431886728379
$ws.Range("P2").Value = 0.9729431886728379
$ws.Range("Q2").Value = 6881.814150815475
$ws.Range("R2").Value = 61936.32735733927
$ws.Range("S2").Value = 0.8103609741962308
$ws.Range("T2").Value = 0.8103609741962309

# Row 3
$ws.Range("G3").Value = 36.94436433333333
$ws.Range("H3").Value = 110.833093
$ws.Range("I3").Value = 0.8328964975864823
$ws.Range("J3").Value = 0.8328964975864824
$ws.Range("O3").Value = 0.01102259370028598
$ws.Range("P3").Value = 0.01102259370028598
$ws.Range("Q3").Value = 77.96492353144444
$ws.Range("R3").Value = 701.684311783
$ws.Range("S3").Value = 0.009180679687287013
$ws.Range("T3").Value = 0.009180679687287014

# Row 4
$ws.Range("G4").Value = 36.94436433333333
$ws.Range("H4").Value = 110.833093
$ws.Range("I4").Value = 0.8328964975864823
$ws.Range("J4").Value = 0.8328964975864824
$ws.Range("M4").Value = 3.069835
$ws.Range("N4").Value = 9.209505
$ws.Range("O4").Value = 0.01603421762687604
$ws.Range("P4").Value = 0.01603421762687604
$ws.Range("Q4").Value = 113.4131026832183
$ws.Range("R4").Value = 1020.717924148965
$ws.Range("S4").Value = 0.01335484370296449
$ws.Range("T4").Value = 0.01335484370296449

# Row 5
$ws.Range("I5").Value = 0.07608399754092349
$ws.Range("J5").Value = 0.07608399754092349
$ws.Range("M5").Value = 186.2750726666667
$ws.Range("N5").Value = 558.825218
$ws.Range("O5").Value = 0.9729431886728379
$ws.Range("P5").Value = 0.9729431886728379
$ws.Range("Q5").Value = 628.6446544618473
$ws.Range("R5").Value = 5657.801890156626
$ws.Range("S5").Value = 0.07402540717444246
$ws.Range("T5").Value = 0.07402540717444246

# Row 6
$ws.Range("I6").Value = 0.07608399754092349
$ws.Range("J6").Value = 0.07608399754092349
$ws.Range("O6").Value = 0.01102259370028598
$ws.Range("P6").Value = 0.01102259370028598
$ws.Range("S6").Value = 0.0008386429919871569
$ws.Range("T6").Value = 0.0008386429919871569

# Row 7
$ws.Range("I7").Value = 0.07608399754092349
$ws.Range("J7").Value = 0.07608399754092349
$ws.Range("M7").Value = 3.069835
$ws.Range("N7").Value = 9.209505
$ws.Range("O7").Value = 0.01603421762687604
$ws.Range("P7").Value = 0.01603421762687604
$ws.Range("Q7").Value = 10.360137484865
$ws.Range("R7").Value = 93.241237363785
$ws.Range("S7").Value = 0.001219947374493869
$ws.Range("T7").Value = 0.001219947374493869

# Row 8
$ws.Range("G8").Value = 4.037305666666668
$ws.Range("H8").Value = 12.111917
$ws.Range("I8").Value = 0.09101950487259411
$ws.Range("J8").Value = 0.09101950487259411
$ws.Range("M8").Value = 186.2750726666667
$ws.Range("N8").Value = 558.825218
$ws.Range("O8").Value = 0.9729431886728379
$ws.Range("P8").Value = 0.9729431886728379
$ws.Range("Q8").Value = 752.0494064358786
$ws.Range("R8").Value = 6768.444657922906
$ws.Range("S8").Value = 0.08855680730216461
$ws.Range("T8").Value = 0.08855680730216461

# Row 9
$ws.Range("G9").Value = 4.037305666666668
$ws.Range("H9").Value = 12.111917
$ws.Range("I9").Value = 0.09101950487259411
$ws.Range("J9").Value = 0.09101950487259411
$ws.Range("O9").Value = 0.01102259370028598
$ws.Range("P9").Value = 0.01102259370028598
$ws.Range("Q9").Value = 8.520060725222223
$ws.Range("R9").Value = 76.680546527
$ws.Range("S9").Value = 0.001003271021011805
$ws.Range("T9").Value = 0.001003271021011805

# Row 10
$ws.Range("G10").Value = 4.037305666666668
$ws.Range("H10").Value = 12.111917
$ws.Range("I10").Value = 0.09101950487259411
$ws.Range("J10").Value = 0.09101950487259411
$ws.Range("M10").Value = 3.069835
$ws.Range("N10").Value = 9.209505
$ws.Range("O10").Value = 0.01603421762687604
$ws.Range("P10").Value = 0.01603421762687604
$ws.Range("Q10").Value = 12.39386224123167
$ws.Range("R10").Value = 111.544760171085
$ws.Range("S10").Value = 0.001459426549417678
$ws.Range("T10").Value = 0.001459426549417678

Write-Output "updated TPM values"
